# Updates cryptos list price/volume columns (inline-string cells).
# A leading apostrophe forces Excel to keep the text verbatim (no numeric
# coercion of values like "211.35"); re-applying the "Normal" style strips
# the quote-prefix formatting Excel would otherwise stamp on the cell, so
# the result matches the original plain (un-styled) inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.721.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.37%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.601.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.36%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'211.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.51%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.27%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.60%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.48%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.64%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.826.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.29%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.601.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.78%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D15").Value = "'0.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.59%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.74%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.692.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0745"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'210.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.19%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.96%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.65%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'8.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.12%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'143.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.27%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.12%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.95%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.69%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.17%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.57%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.93%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.298.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.05%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.66%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.86%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.16%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +20.74%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.27%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.76%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.27%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.15%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.782"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.43%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'63.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.13%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.737.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'91.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.24%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.99%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.41%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.64%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +1.98%  "
$ws.Range("E50").Style = "Normal"
